$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Predictions")

# Row 3 - San Francisco 49ers @ Jacksonville Jaguars
$ws.Range("D3").Value = 13
$ws.Range("E3").Value = 14
$ws.Range("H3").Value = 40
$ws.Range("J3").Value = "UNDER"

# Row 4 - Indianapolis Colts @ Buffalo Bills
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = 42
$ws.Range("E4").Value = -20
$ws.Range("H4").Value = 64

# Row 5 - Green Bay Packers @ Minnesota Vikings
$ws.Range("C5").Value = 21
$ws.Range("D5").Value = 23
$ws.Range("E5").Value = -2
$ws.Range("G5").Value = "Minnesota +2"
$ws.Range("H5").Value = 42

# Row 6 - Baltimore Ravens @ Chicago Bears
$ws.Range("C6").Value = 31
$ws.Range("D6").Value = 13
$ws.Range("E6").Value = 18
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = "Baltimore -5"
$ws.Range("H6").Value = 44
$ws.Range("J6").Value = "UNDER"

# Row 7 - Washington Football Team @ Carolina Panthers
$ws.Range("C7").Value = 14
$ws.Range("D7").Value = 25
$ws.Range("G7").Value = "Carolina -3.5"
$ws.Range("H7").Value = 39
$ws.Range("J7").Value = "UNDER"

# Row 8 - Detroit Lions @ Cleveland Browns
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = 31
$ws.Range("E8").Value = 20
$ws.Range("G8").Value = "Cleveland -11.5"
$ws.Range("H8").Value = 42
$ws.Range("J8").Value = "UNDER"

# Row 9 - Houston Texans @ Tennessee Titans
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 39
$ws.Range("E9").Value = -33
$ws.Range("G9").Value = "Tennessee - 10"
$ws.Range("H9").Value = 45

# Row 10 - New Orleans Saints @ Philadelphia Eagles
$ws.Range("C10").Value = 26
$ws.Range("D10").Value = 26
$ws.Range("G10").Value = "New Orleans +1.5"
$ws.Range("H10").Value = 52

# Row 11 - Miami Dolphins @ New York Jets
$ws.Range("C11").Value = 21
$ws.Range("D11").Value = 18
$ws.Range("E11").Value = 3
$ws.Range("H11").Value = 39
$ws.Range("J11").Value = "UNDER"

# Row 12 - Cincinnati Bengals @ Las Vegas Raiders
$ws.Range("C12").Value = 32
$ws.Range("D12").Value = 27
$ws.Range("E12").Value = 5
$ws.Range("H12").Value = 59

# Row 13 - Dallas Cowboys @ Kansas City Chiefs
$ws.Range("C13").Value = 39
$ws.Range("D13").Value = 31
$ws.Range("E13").Value = 8
$ws.Range("G13").Value = "Dallas -2.5"
$ws.Range("H13").Value = 70

# Row 14 - Arizona Cardinals @ Seattle Seahawks
$ws.Range("G14").Value = "Waiting on Injury Report"

# Row 15 - Pittsburgh Steelers @ Los Angeles Chargers
$ws.Range("C15").Value = 18
$ws.Range("D15").Value = 27
$ws.Range("E15").Value = -9
$ws.Range("G15").Value = "Los Angeles Chargers -5.5"
$ws.Range("H15").Value = 45
$ws.Range("J15").Value = "UNDER"

# Row 16 - New York Giants @ Tampa Bay Buccaneers
$ws.Range("C16").Value = 16
$ws.Range("D16").Value = 42
$ws.Range("E16").Value = -26
$ws.Range("G16").Value = "Tampa Bay -11.5"
$ws.Range("H16").Value = 58

# Update the view state: activate Predictions sheet, scroll and select D8
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D8").Select()
